$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 330
$ws1.Range("F9").Value = 13024
$ws1.Range("F12").Value = 5288
$ws1.Range("F21").Value = 2863
$ws1.Range("F22").Value = 6215

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 330
$ws4.Range("F10").Value = 13024
$ws4.Range("F13").Value = 5288
$ws4.Range("F22").Value = 2863
$ws4.Range("F24").Value = 6215
